$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.857000000000001
$ws.Range("D4").Value = -8.074999999999999
$ws.Range("D7").Value = -7.935
$ws.Range("D8").Value = -8.108000000000001
$ws.Range("B11").Value = 6.661
$ws.Range("B12").Value = 5.915
$ws.Range("D12").Value = -8.311
$ws.Range("D14").Value = -8.178999999999998
$ws.Range("B15").Value = 6.441
$ws.Range("D22").Value = -7.812
